$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.899.00'
$ws.Range("E2").Value = '  -0.67%  '
$ws.Range("D3").Value = '1.667.88'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("E6").Value = '  +5.49%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +1.09%  '
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.25'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0896'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.63%  '
$ws.Range("D12").Value = '1.903.07'
$ws.Range("E12").Value = '  +0.66%  '
$ws.Range("D13").Value = '1.678.65'
$ws.Range("E13").Value = '  +1.34%  '
$ws.Range("E14").Value = '  +0.06%  '
$ws.Range("E15").Value = '  +1.18%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.19'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.52%  '
$ws.Range("D17").Value = '26.913.75'
$ws.Range("E17").Value = '  -0.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '234.43'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.95'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.23%  '
$ws.Range("D20").Value = '0.0₃0732'
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.72%  '
$ws.Range("B23").Value = 'Toncoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.97%  '
$ws.Range("B24").Value = 'Avalanche'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.14'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.31'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.26%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  +0.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.89'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.26%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("E32").Value = '  +1.76%  '
$ws.Range("D33").Value = '1.452.89'
$ws.Range("E33").Value = '  -4.25%  '
$ws.Range("E34").Value = '  +2.18%  '
$ws.Range("E35").Value = '  +2.71%  '
$ws.Range("E36").Value = '  -0.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.581'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.905'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.87%  '
$ws.Range("E39").Value = '  +0.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.73'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.80%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("E42").Value = '  +0.83%  '
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.974'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.93%  '
$ws.Range("D45").Value = '1.808.95'
$ws.Range("E45").Value = '  +0.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.784'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.07%  '
$ws.Range("E48").Value = '  +1.04%  '
$ws.Range("E49").Value = '  -1.11%  '
$ws.Range("E50").Value = '  +4.49%  '
$ws.Range("E51").Value = '  -0.13%  '
